# Apply ROI parameter changes as described in the commit:
# - Row 5 ("Bec") measured ROI coordinates & calibration value were updated.
# - A new row 9 ("NiLattice") was added with its own ROI coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (Bec) ---
$ws.Range("A5").Value = "Bec"
$ws.Range("B5").Value = 1172
$ws.Range("C5").Value = 1350
$ws.Range("D5").Value = 1667
$ws.Range("E5").Value = 1881
$ws.Range("F5").Value = 2160
$ws.Range("G5").Value = 2560
$ws.Range("H5").Value = 349.10000000000002

# --- Add new row 9 (NiLattice) ---
$ws.Range("A9").Value = "NiLattice"
$ws.Range("B9").Value = 1107
$ws.Range("C9").Value = 1417
$ws.Range("D9").Value = 1740
$ws.Range("E9").Value = 1814
$ws.Range("F9").Value = 2160
$ws.Range("G9").Value = 2560
$ws.Range("H9").Value = 349.10000000000002
